$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.965.45"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").Value = "1.859.76"
$ws.Range("E3").Value = "  -2.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4361"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.85%  "
$ws.Range("E8").Value = "  -3.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9408"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.65%  "
$ws.Range("D12").Value = "1.863.94"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06844"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009035"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.20%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.90%  "
$ws.Range("D21").Value = "27.938.51"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("E22").Value = "  -3.79%  "
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").Value = "2.094.27"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.005"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.392"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.737"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.22%  "
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8115"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.88%  "
$ws.Range("E33").Value = "  -5.33%  "
$ws.Range("E34").Value = "  -5.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.932"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.68%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05493"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.113"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01973"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.42%  "
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5260"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.018"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.48%  "
$ws.Range("E43").Value = "  -3.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.795"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06790"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4889"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.77%  "
$ws.Range("E49").Value = "  -5.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.913"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -10.26%  "
$ws.Range("E51").Value = "  -0.14%  "
